$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.031884919460599
$ws.Range("D2").Value = 1.039780398609737
$ws.Range("E2").Value = 1.040345215763777
$ws.Range("F2").Value = 1.048735872445487
$ws.Range("I2").Value = 1.02592303400647
$ws.Range("J2").Value = 1.037018184258537
$ws.Range("K2").Value = 1.042564613932337
$ws.Range("L2").Value = 1.043127827785265
$ws.Range("M2").Value = 1.051494887174389
$ws.Range("N2").Value = 1.016327136406408

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.033087756228901
$ws.Range("D3").Value = 1.040904770533465
$ws.Range("E3").Value = 1.041463438377412
$ws.Range("F3").Value = 1.050068373434721
$ws.Range("I3").Value = 1.026089590951505
$ws.Range("J3").Value = 1.037861795254908
$ws.Range("K3").Value = 1.043498447266986
$ws.Range("L3").Value = 1.044055644888622
$ws.Range("M3").Value = 1.052638148005665
$ws.Range("N3").Value = 1.016605090389138

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.033860997972303
$ws.Range("D4").Value = 1.041626102621187
$ws.Range("E4").Value = 1.042180757432614
$ws.Range("F4").Value = 1.050919091200678
$ws.Range("I4").Value = 1.026187754699131
$ws.Range("J4").Value = 1.03840212569305
$ws.Range("K4").Value = 1.044095922251125
$ws.Range("L4").Value = 1.044649191463329
$ws.Range("M4").Value = 1.053365903434643
$ws.Range("N4").Value = 1.016783110577857

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.034184865296069
$ws.Range("D5").Value = 1.041927875127717
$ws.Range("E5").Value = 1.042480834711244
$ws.Range("F5").Value = 1.051273997879461
$ws.Range("I5").Value = 1.026226726306112
$ws.Range("J5").Value = 1.038627963114772
$ws.Range("K5").Value = 1.044345489014111
$ws.Range("L5").Value = 1.044897098224083
$ws.Range("M5").Value = 1.053668992874885
$ws.Range("N5").Value = 1.01685751413139

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.034239173784485
$ws.Range("D6").Value = 1.041978457856984
$ws.Range("E6").Value = 1.0425311323262
$ws.Range("F6").Value = 1.051333428382213
$ws.Range("I6").Value = 1.026233135251464
$ws.Range("J6").Value = 1.038665805211978
$ws.Range("K6").Value = 1.044387298160049
$ws.Range("L6").Value = 1.04493862814992
$ws.Range("M6").Value = 1.053719715728934
$ws.Range("N6").Value = 1.01686998132571

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.033865330217269
$ws.Range("D7").Value = 1.04163014070249
$ws.Range("E7").Value = 1.042184772894229
$ws.Range("F7").Value = 1.050923844204359
$ws.Range("I7").Value = 1.026188284457674
$ws.Range("J7").Value = 1.038405148508178
$ws.Range("K7").Value = 1.044099263292989
$ws.Range("L7").Value = 1.044652510357559
$ws.Range("M7").Value = 1.053369964542553
$ws.Range("N7").Value = 1.016784106471396

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.032292481392361
$ws.Range("D8").Value = 1.040161678795325
$ws.Range("E8").Value = 1.040724424802958
$ws.Range("F8").Value = 1.049188589671269
$ws.Range("I8").Value = 1.025981315236096
$ws.Range("J8").Value = 1.037304439954335
$ws.Range("K8").Value = 1.042881617433399
$ws.Range("L8").Value = 1.043442805426372
$ws.Range("M8").Value = 1.051883754784169
$ws.Range("N8").Value = 1.01642145405254

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.02948151765874
$ws.Range("D9").Value = 1.037525921935185
$ws.Range("E9").Value = 1.038102705002597
$ws.Range("F9").Value = 1.04604189900554
$ws.Range("I9").Value = 1.025542776820477
$ws.Range("J9").Value = 1.035321947860525
$ws.Range("K9").Value = 1.040683529981103
$ws.Range("L9").Value = 1.041258438708941
$ws.Range("M9").Value = 1.049172063069363
$ws.Range("N9").Value = 1.015768213046998

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.027580236189411
$ws.Range("D10").Value = 1.035735528660811
$ws.Range("E10").Value = 1.036321494317958
$ws.Range("F10").Value = 1.043882963226861
$ws.Range("I10").Value = 1.02520038170333
$ws.Range("J10").Value = 1.033970749572417
$ws.Range("K10").Value = 1.039182072426321
$ws.Range("L10").Value = 1.039765944624168
$ws.Range("M10").Value = 1.047300631330381
$ws.Range("N10").Value = 1.015322943510497

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.026750286924182
$ws.Range("D11").Value = 1.034952179508049
$ws.Range("E11").Value = 1.035542079560836
$ws.Range("F11").Value = 1.04293328967074
$ws.Range("I11").Value = 1.02504014511123
$ws.Range("J11").Value = 1.033378490735136
$ws.Range("K11").Value = 1.038523171859824
$ws.Range("L11").Value = 1.039110881534008
$ws.Range("M11").Value = 1.04647487924513
$ws.Range("N11").Value = 1.015127762590115

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.026440984897
$ws.Range("D12").Value = 1.034659972950996
$ws.Range("E12").Value = 1.035251328125568
$ws.Range("F12").Value = 1.042578278655305
$ws.Range("I12").Value = 1.024978815975375
$ws.Range("J12").Value = 1.033157404915053
$ws.Range("K12").Value = 1.038277092452249
$ws.Range("L12").Value = 1.038866220934579
$ws.Range("M12").Value = 1.046165814397414
$ws.Range("N12").Value = 1.015054901556089

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.026507377775784
$ws.Range("D13").Value = 1.034722708494214
$ws.Range("E13").Value = 1.035313751825119
$ws.Range("F13").Value = 1.04265453249972
$ws.Range("I13").Value = 1.024992053377998
$ws.Range("J13").Value = 1.0332048783642
$ws.Range("K13").Value = 1.038329938011782
$ws.Range("L13").Value = 1.038918762454145
$ws.Range("M13").Value = 1.046232216394952
$ws.Range("N13").Value = 1.015070546969028

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.026724740876676
$ws.Range("D14").Value = 1.034928050948538
$ws.Range("E14").Value = 1.035518071410544
$ws.Range("F14").Value = 1.042903990634568
$ws.Range("I14").Value = 1.025035112617277
$ws.Range("J14").Value = 1.033360238133812
$ws.Range("K14").Value = 1.038502858190741
$ws.Range("L14").Value = 1.039090685278883
$ws.Range("M14").Value = 1.046449379832137
$ws.Range("N14").Value = 1.01512174729285

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.02685852948759
$ws.Range("D15").Value = 1.035054404859826
$ws.Range("E15").Value = 1.035643794264686
$ws.Range("F15").Value = 1.043057389616446
$ws.Range("I15").Value = 1.02506140263912
$ws.Range("J15").Value = 1.033455814898192
$ws.Range("K15").Value = 1.038609222631388
$ws.Range("L15").Value = 1.039196434341065
$ws.Range("M15").Value = 1.046582869937195
$ws.Range("N15").Value = 1.015153245355113

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.027635174890876
$ws.Range("D16").Value = 1.03578734474474
$ws.Range("E16").Value = 1.036373048408357
$ws.Range("F16").Value = 1.043945674713279
$ws.Range("I16").Value = 1.025210762842183
$ws.Range("J16").Value = 1.034009903228937
$ws.Range("K16").Value = 1.039225615388944
$ws.Range("L16").Value = 1.039809231981907
$ws.Range("M16").Value = 1.047355106707667
$ws.Range("N16").Value = 1.015335846527726

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.028120542666399
$ws.Range("D17").Value = 1.036244917464832
$ws.Range("E17").Value = 1.036828297976999
$ws.Range("F17").Value = 1.044498878451169
$ws.Range("I17").Value = 1.025301238674405
$ws.Range("J17").Value = 1.034355534858682
$ws.Range("K17").Value = 1.039609904462531
$ws.Range("L17").Value = 1.040191254199464
$ws.Range("M17").Value = 1.047835365140952
$ws.Range("N17").Value = 1.015449747684028

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.028403006048594
$ws.Range("D18").Value = 1.036511032008794
$ws.Range("E18").Value = 1.037093053441584
$ws.Range("F18").Value = 1.044820123130535
$ws.Range("I18").Value = 1.025352856846279
$ws.Range("J18").Value = 1.034556444179192
$ws.Range("K18").Value = 1.039833209910808
$ws.Range("L18").Value = 1.040413233208276
$ws.Range("M18").Value = 1.048114006969506
$ws.Range("N18").Value = 1.015515955389251

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.028499210156794
$ws.Range("D19").Value = 1.036601638476841
$ws.Range("E19").Value = 1.037183195816576
$ws.Range("F19").Value = 1.044929417654346
$ws.Range("I19").Value = 1.025370261733321
$ws.Range("J19").Value = 1.034624832190018
$ws.Range("K19").Value = 1.03990920870265
$ws.Range("L19").Value = 1.040488779003357
$ws.Range("M19").Value = 1.048208765624032
$ws.Range("N19").Value = 1.01553849182926

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.028068533978368
$ws.Range("D20").Value = 1.036195905045823
$ws.Range("E20").Value = 1.036779535226545
$ws.Range("F20").Value = 1.04443967302105
$ws.Range("I20").Value = 1.025291650999754
$ws.Range("J20").Value = 1.034318523518563
$ws.Range("K20").Value = 1.039568761296737
$ws.Range("L20").Value = 1.040150354682609
$ws.Range("M20").Value = 1.047783991736799
$ws.Range("N20").Value = 1.015437550881722

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.026660761220464
$ws.Range("D21").Value = 1.034867616988372
$ws.Range("E21").Value = 1.035457938833527
$ws.Range("F21").Value = 1.042830594070659
$ws.Range("I21").Value = 1.02502248280375
$ws.Range("J21").Value = 1.033314518897911
$ws.Range("K21").Value = 1.038451974496713
$ws.Range("L21").Value = 1.039040095463278
$ws.Range("M21").Value = 1.046385495549979
$ws.Range("N21").Value = 1.01510668011343

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.025769716515237
$ws.Range("D22").Value = 1.034025309250928
$ws.Range("E22").Value = 1.034619802018423
$ws.Range("F22").Value = 1.041805810432374
$ws.Range("I22").Value = 1.024842767204449
$ws.Range("J22").Value = 1.032676920420429
$ws.Range("K22").Value = 1.037742076021844
$ws.Range("L22").Value = 1.038334262927466
$ws.Range("M22").Value = 1.045492629820329
$ws.Range("N22").Value = 1.014896550488686

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.026242643784703
$ws.Range("D23").Value = 1.034472517963035
$ws.Range("E23").Value = 1.035064803102803
$ws.Range("F23").Value = 1.042350319197268
$ws.Range("I23").Value = 1.024939034947364
$ws.Range("J23").Value = 1.033015529935817
$ws.Range("K23").Value = 1.038119145788395
$ws.Range("L23").Value = 1.038709180894566
$ws.Range("M23").Value = 1.045967251657155
$ws.Range("N23").Value = 1.015008144840057

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.028092036453504
$ws.Range("D24").Value = 1.036218054057502
$ws.Range("E24").Value = 1.036801571436182
$ws.Range("F24").Value = 1.044466429823109
$ws.Range("I24").Value = 1.025295986825361
$ws.Range("J24").Value = 1.034335249487071
$ws.Range("K24").Value = 1.03958735472931
$ws.Range("L24").Value = 1.040168838033752
$ws.Range("M24").Value = 1.047807209753721
$ws.Range("N24").Value = 1.015443062798854

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.030212969209207
$ws.Range("D25").Value = 1.038213117447244
$ws.Range("E25").Value = 1.038786302485109
$ws.Range("F25").Value = 1.046866066599899
$ws.Range("I25").Value = 1.025664936159382
$ws.Range("J25").Value = 1.035839621964143
$ws.Range("K25").Value = 1.041258082250312
$ws.Range("L25").Value = 1.041829476149565
$ws.Range("M25").Value = 1.049884219724321
$ws.Range("N25").Value = 1.015938796845746
